$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1692307692307692
$ws.Range("C2").Value = 0.5815384615384616
$ws.Range("J2").Value = 0.02153846153846154
$ws.Range("P2").Value = 0.1353846153846154
$ws.Range("S2").Value = 0.09230769230769231
$ws.Range("B3").Value = 0.005208333333333333
$ws.Range("C3").Value = 0.01041666666666667
$ws.Range("J3").Value = 0.02604166666666667
$ws.Range("P3").Value = 0.7291666666666666
$ws.Range("S3").Value = 0.2291666666666667
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.7115384615384616
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.0743801652892562
$ws.Range("D6").Value = 0.01239669421487603
$ws.Range("F6").Value = 0.06611570247933884
$ws.Range("J6").Value = 0.231404958677686
$ws.Range("O6").Value = 0.04132231404958678
$ws.Range("Q6").Value = 0.1611570247933884
$ws.Range("R6").Value = 0.07851239669421488
$ws.Range("S6").Value = 0.3347107438016529
$ws.Range("B7").Value = 0.1152073732718894
$ws.Range("F7").Value = 0.03225806451612903
$ws.Range("J7").Value = 0.152073732718894
$ws.Range("O7").Value = 0.04608294930875576
$ws.Range("Q7").Value = 0.184331797235023
$ws.Range("R7").Value = 0.08755760368663594
$ws.Range("S7").Value = 0.3824884792626728
$ws.Range("B8").Value = 0.09430255402750491
$ws.Range("D8").Value = 0.0137524557956778
$ws.Range("E8").Value = 0.003929273084479371
$ws.Range("F8").Value = 0.07269155206286837
$ws.Range("J8").Value = 0.1139489194499018
$ws.Range("O8").Value = 0.009823182711198428
$ws.Range("Q8").Value = 0.1611001964636542
$ws.Range("R8").Value = 0.1296660117878193
$ws.Range("S8").Value = 0.4007858546168959
$ws.Range("B9").Value = 0.06617647058823529
$ws.Range("D9").Value = 0.03676470588235294
$ws.Range("F9").Value = 0.06617647058823529
$ws.Range("J9").Value = 0.1029411764705882
$ws.Range("O9").Value = 0.02205882352941177
$ws.Range("Q9").Value = 0.1176470588235294
$ws.Range("R9").Value = 0.1911764705882353
$ws.Range("S9").Value = 0.3970588235294117
$ws.Range("B10").Value = 0.1237541528239203
$ws.Range("D10").Value = 0.03239202657807309
$ws.Range("E10").Value = 0.0008305647840531562
$ws.Range("F10").Value = 0.06810631229235881
$ws.Range("J10").Value = 0.117109634551495
$ws.Range("O10").Value = 0.02491694352159468
$ws.Range("Q10").Value = 0.1794019933554817
$ws.Range("R10").Value = 0.1096345514950166
$ws.Range("S10").Value = 0.3438538205980066
$ws.Range("G11").Value = 0.1596638655462185
$ws.Range("J11").Value = 0.08683473389355742
$ws.Range("K11").Value = 0.2212885154061625
$ws.Range("L11").Value = 0.5098039215686274
$ws.Range("S11").Value = 0.02240896358543417
$ws.Range("G12").Value = 0.7540983606557377
$ws.Range("J12").Value = 0.2021857923497268
$ws.Range("K12").Value = 0.01639344262295082
$ws.Range("L12").Value = 0.01092896174863388
$ws.Range("S12").Value = 0.01639344262295082
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.1666666666666667
$ws.Range("S13").Value = 0.119047619047619
$ws.Range("F15").Value = 0.008733624454148471
$ws.Range("H15").Value = 0.2008733624454148
$ws.Range("I15").Value = 0.03930131004366812
$ws.Range("J15").Value = 0.3362445414847162
$ws.Range("K15").Value = 0.08733624454148471
$ws.Range("M15").Value = 0.01746724890829694
$ws.Range("N15").Value = 0.004366812227074236
$ws.Range("O15").Value = 0.07860262008733625
$ws.Range("S15").Value = 0.2270742358078603
$ws.Range("F16").Value = 0.03381642512077294
$ws.Range("H16").Value = 0.1980676328502415
$ws.Range("I16").Value = 0.05314009661835749
$ws.Range("J16").Value = 0.391304347826087
$ws.Range("K16").Value = 0.0821256038647343
$ws.Range("M16").Value = 0.02415458937198068
$ws.Range("O16").Value = 0.07246376811594203
$ws.Range("S16").Value = 0.1449275362318841
$ws.Range("F17").Value = 0.02798982188295165
$ws.Range("H17").Value = 0.2010178117048346
$ws.Range("I17").Value = 0.04834605597964377
$ws.Range("J17").Value = 0.3969465648854962
$ws.Range("K17").Value = 0.1374045801526718
$ws.Range("M17").Value = 0.02035623409669211
$ws.Range("O17").Value = 0.05089058524173028
$ws.Range("S17").Value = 0.1170483460559796
$ws.Range("F18").Value = 0.02298850574712644
$ws.Range("H18").Value = 0.1877394636015326
$ws.Range("I18").Value = 0.09578544061302682
$ws.Range("J18").Value = 0.4367816091954023
$ws.Range("K18").Value = 0.08812260536398467
$ws.Range("M18").Value = 0.01149425287356322
$ws.Range("N18").Value = 0.003831417624521073
$ws.Range("O18").Value = 0.05363984674329502
$ws.Range("S18").Value = 0.09961685823754789
$ws.Range("F19").Value = 0.0322061191626409
$ws.Range("H19").Value = 0.2375201288244767
$ws.Range("I19").Value = 0.0571658615136876
$ws.Range("J19").Value = 0.3276972624798712
$ws.Range("K19").Value = 0.1191626409017713
$ws.Range("M19").Value = 0.01932367149758454
$ws.Range("O19").Value = 0.06521739130434782
$ws.Range("S19").Value = 0.14170692431562
